# Update "Forecast Comparison" sheet:
#  - insert a new column B "Week_Start_Date" (ASIN and everything to the right
#    shifts over by one column; no other columns are added/removed)
#  - shorten the week labels in column A (W01 -> W1 ... W09 -> W9; W10+ unchanged)
#  - fill the new Week_Start_Date column with the Monday date of each week
#  - is_holiday_week (now column J) becomes a boolean value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert new column B (Week_Start_Date); everything from ASIN onward shifts right ---
$ws.Range("B:B").Insert()

# --- Header ---
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week_Start_Date values are stored as plain text (e.g. "2025-01-05"), not
# as Excel date serials, so force the column to text format before writing.
$ws.Range("B2:B17").NumberFormat = "@"

# --- Week start dates (Monday) for each of the 16 weekly rows ---
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    # Shorten the week label (W01 -> W1 ... W09 -> W9; W10+ stays as-is)
    $week = $ws.Cells.Item($row, 1).Value2
    if ($week -match '^W0(\d)$') {
        $ws.Cells.Item($row, 1).Value = "W" + $matches[1]
    }

    # New Week_Start_Date column
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]

    # is_holiday_week (now column J) stored as boolean
    $ws.Cells.Item($row, 10).Value = $false
}
